# Apply crypto price/volume updates from the Oct 12 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.763.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.38%  "

$ws.Range("D3").Value = "'2.445.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.01%  "

$ws.Range("D5").Value = "'576.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.69%  "

$ws.Range("D6").Value = "'145.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.97%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").Value = "'2.445.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("E10").Value = "  +3.13%  "

$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("D13").Value = "'0.353"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.48%  "

$ws.Range("D14").Value = "'28.17"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.72%  "

$ws.Range("D15").Value = "'0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.84%  "

$ws.Range("E16").Value = "  +2.75%  "

$ws.Range("D17").Value = "'62.629.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.74%  "

$ws.Range("D18").Value = "'2.447.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.79%  "

$ws.Range("D19").Value = "'7.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.75%  "

$ws.Range("E20").Value = "  +2.85%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'330.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("B22").Value = "BabyDogeCoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D22").Value = "'0.0₆0785"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +182.00%  "

$ws.Range("E24").Value = "  +9.43%  "

$ws.Range("D26").Value = "'65.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.05%  "

$ws.Range("D27").Value = "'645.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +13.08%  "

$ws.Range("E28").Value = "  +17.38%  "

$ws.Range("D29").Value = "'8.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.76%  "

$ws.Range("D30").Value = "'0.0₃0986"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.11%  "

$ws.Range("D32").Value = "'1.45"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.74%  "

$ws.Range("E33").Value = "  +1.80%  "

$ws.Range("E34").Value = "  +3.76%  "

$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("E36").Value = "  +2.66%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  +3.43%  "

$ws.Range("D39").Value = "'5.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.21%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.374"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.01%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'153.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.75%  "

$ws.Range("D42").Value = "'18.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.53%  "

$ws.Range("D43").Value = "'2.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.72%  "

$ws.Range("E44").Value = "  +4.83%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'42.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.50%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "'14.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +27.30%  "

$ws.Range("D48").Value = "'145.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.38%  "

$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("D50").Value = "'20.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.45%  "

$ws.Range("D51").Value = "'0.604"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.78%  "
